$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-list row was reported for the week of 2020-11-24 (Rainier / Especial).
# Insert a fresh row at position 289 (pushes every following row down by one,
# which also re-aligns the trailing rows/dates exactly as they were before,
# just shifted), then populate the new row with its data.
$ws.Rows.Item(289).Insert()

$ws.Cells.Item(289, 1).Value = 6
$ws.Cells.Item(289, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(289, 3).Value = 'Metropolitana'
$ws.Cells.Item(289, 4).Value = 44159
$ws.Cells.Item(289, 5).Value = 13
$ws.Cells.Item(289, 6).Value = 'Fruta'
$ws.Cells.Item(289, 7).Value = 100103
$ws.Cells.Item(289, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(289, 9).Value = 100103001
$ws.Cells.Item(289, 10).Value = 'Cereza'
$ws.Cells.Item(289, 11).Value = 'Rainier'
$ws.Cells.Item(289, 12).Value = 'Especial'
$ws.Cells.Item(289, 13).Value = 200
$ws.Cells.Item(289, 14).Value = 23000
$ws.Cells.Item(289, 15).Value = 23000
$ws.Cells.Item(289, 16).Value = 23000
$ws.Cells.Item(289, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(289, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(289, 19).Value = 2300
$ws.Cells.Item(289, 20).Value = 10
